$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------
# Helper data: the four paragraphs that become bold "section
# headers" for the reference list, identified by their exact text.
# ---------------------------------------------------------------

# --- 1) "Original paper describing the design:" ---
$p = $d.Paragraphs.Item(3)
$orig = $p.Range.WordOpenXML
$null = $orig -match '(?s)<pkg:part pkg:name="/word/document.xml"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData></pkg:part>'
$docXml = $matches[1]
$marker = "Original paper describing the design:"
$null = $docXml -match ('(?s)<w:p(\s[^>]*)?>.*?' + [regex]::Escape($marker))
$attrs = $matches[1]
$newXml = $pkgHeader + '<w:p' + $attrs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>' + $marker + '</w:t></w:r></w:p>' + $pkgFooter
$p.Range.InsertXML($newXml)

# --- 2) "Literature review of use of the PNUD:" ---
$p = $d.Paragraphs.Item(5)
$orig = $p.Range.WordOpenXML
$null = $orig -match '(?s)<pkg:part pkg:name="/word/document.xml"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData></pkg:part>'
$docXml = $matches[1]
$marker = "Literature review of use of the PNUD:"
$null = $docXml -match ('(?s)<w:p(\s[^>]*)?>.*?' + [regex]::Escape($marker))
$attrs = $matches[1]
$newXml = $pkgHeader + '<w:p' + $attrs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>' + $marker + '</w:t></w:r></w:p>' + $pkgFooter
$p.Range.InsertXML($newXml)

# --- 3) "Simulation study that showed ... (in the scenario of a single switch):" ---
$p = $d.Paragraphs.Item(7)
$orig = $p.Range.WordOpenXML
$null = $orig -match '(?s)<pkg:part pkg:name="/word/document.xml"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData></pkg:part>'
$docXml = $matches[1]
$marker = "Simulation study that showed that conditioning on time alone can result in a biased estimate (in the scenario of a single switch):"
$null = $docXml -match ('(?s)<w:p(\s[^>]*)?>.*?' + [regex]::Escape($marker))
$attrs = $matches[1]
$newXml = $pkgHeader + '<w:p' + $attrs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>' + $marker + '</w:t></w:r></w:p>' + $pkgFooter
$p.Range.InsertXML($newXml)

# --- 4) "Example of use of the PNUD for patients with more complex treatment histories: " ---
$p = $d.Paragraphs.Item(9)
$orig = $p.Range.WordOpenXML
$null = $orig -match '(?s)<pkg:part pkg:name="/word/document.xml"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData></pkg:part>'
$docXml = $matches[1]
$marker = "Example of use of the PNUD for patients with more complex treatment histories: "
$null = $docXml -match ('(?s)<w:p(\s[^>]*)?>.*?' + [regex]::Escape($marker))
$attrs = $matches[1]
$newXml = $pkgHeader + '<w:p' + $attrs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">' + $marker + '</w:t></w:r></w:p>' + $pkgFooter
$p.Range.InsertXML($newXml)

# ---------------------------------------------------------------
# Insert the two new paragraphs (new bold header + new citation)
# right after the Webster-Clark citation (paragraph 8) and before
# "Example of use of the PNUD ..." (paragraph 9).
# ---------------------------------------------------------------

$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()

$newHeading = $d.Paragraphs.Item(9)
$headingText = "Description of the different types of new-user studies used in pharmacoepidemiology:"
$headingXml = $pkgHeader + '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>' + $headingText + '</w:t></w:r></w:p>' + $pkgFooter
$newHeading.Range.InsertXML($headingXml)

$newHeading = $d.Paragraphs.Item(9)
$newHeading.Range.InsertParagraphAfter()

$newCitation = $d.Paragraphs.Item(10)
$citationXml = $pkgHeader + '<w:p><w:r><w:t xml:space="preserve">Her QL, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rouette</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> J, Young JC, Webster-Clark M, Tazare J. Core Concepts in Pharmacoepidemiology: New-User Designs. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pharmacoepidemiol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Drug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Saf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. 2024 Dec;33(12</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>):e</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">70048. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: 10.1002/pds.70048.</w:t></w:r></w:p>' + $pkgFooter
$newCitation.Range.InsertXML($citationXml)

Write-Output "All edits applied."
